# Update the "Team" column (4th column) of the Presentation Schedule table.
# Each row's 1st column identifies the week; map week -> assigned team
# per the commit message / diff.

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$assignments = @{
    "Week 4"  = "Team 10"
    "Week 5"  = "Team 8"
    "Week 6"  = "Team 3"
    "Week 7"  = "Team 1"
    "Week 9"  = "Team 6"
    "Week 10" = "Team 5"
    "Week 11" = "Team 11"
    "Week 12" = "Team 4"
    "Week 13" = "Team 9"
    "Week 14" = "Team 7"
    "Week 15" = "Team 2"
}

for ($i = 1; $i -le $t.Rows.Count; $i++) {
    $weekCell = $t.Cell($i, 1).Range.Text
    $week = $weekCell.TrimEnd([char]7, [char]13, [char]10)

    if ($assignments.ContainsKey($week)) {
        $teamCell = $t.Cell($i, 4)
        $cellRange = $teamCell.Range
        $cellRange.End = $cellRange.End - 1
        $cellRange.Text = $assignments[$week]
    }
}
